# Add support for reports grouped by organ:
# Insert two new columns (Organ ID, Organ Description) between
# "Name" (A) and "Link" (B), so the final layout is:
#   A: Name | B: Organ ID | C: Organ Description | D: Link

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the hyperlink target before the layout changes underneath it.
$linkUrl = "https://tacc.utexas.edu/"

# Insert two columns at B so existing column B ("Link") shifts to D.
$ws.Range("B:C").Insert()

# Header row (inherits the bold header formatting from the insert).
$ws.Range("B1").Value = "Organ ID"
$ws.Range("C1").Value = "Organ Description"

# Data row
$ws.Range("B2").Value = "S1"
$ws.Range("C2").Value = "Sphere 1"

# The column insert leaves the worksheet's hyperlink still anchored on the
# old B2 location even though the cell contents moved to D2 - rebuild it
# pointing at the new location.
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), $linkUrl)
$ws.Range("D2").Style = "Hyperlink"

# Column widths for the two new columns (auto-fit to their contents).
$ws.Columns.Item(2).ColumnWidth = 7.42
$ws.Columns.Item(3).ColumnWidth = 15.42

# Selection, to match the final view state.
$ws.Range("C1").Select()
